# Fruta / hortaliza, semanal
# Insert a new weekly price-report row at the top of the Chirimoya data
# block (row 199), pushing the existing rows 199-219 down to 200-220.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new record: shift rows 199..219 down to 200..220.
$ws.Rows.Item(199).Insert()

# Populate the newly inserted row 199 with this week's data point.
$ws.Range("A199").Value = 10
$ws.Range("B199").Value = "Vega Modelo de Temuco"
$ws.Range("C199").Value = "La Araucanía"
$ws.Range("D199").Value = 45194
$ws.Range("E199").Value = 9
$ws.Range("F199").Value = "Fruta"
$ws.Range("G199").Value = 100107
$ws.Range("H199").Value = "Otros"
$ws.Range("I199").Value = 100107002
$ws.Range("J199").Value = "Chirimoya"
$ws.Range("K199").Value = "Cultivar IV Región"
$ws.Range("L199").Value = "Primera"
$ws.Range("M199").Value = 80
$ws.Range("N199").Value = 2500
$ws.Range("O199").Value = 2500
$ws.Range("P199").Value = 2500
$ws.Range("Q199").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R199").Value = "Provincia del Elquí"
$ws.Range("S199").Value = 2500
$ws.Range("T199").Value = 1
